$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$values = @{
  "H116" = 1928.8928
  "I116" = 1785.7059
  "J116" = 2150.182
  "K116" = 1785.7059
  "L116" = 2150.182
  "M116" = 1656.2941
  "N116" = -9034.182000000001
}
foreach ($addr in $values.Keys) {
  $ws.Range($addr).Value = $values[$addr]
}
$ws.Range("H125:N141").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$values = @{
  "H2" = 1480.6842
  "I2" = 1646
  "J2" = 860.75
  "K2" = 1646
  "L2" = 860.75
  "M2" = -1533
  "N2" = -1086.75
  "H61" = 7582.74
  "I61" = 4498.407
  "J61" = 11203.479
  "K61" = 4498.407
  "L61" = 11203.479
  "M61" = -4286.407
  "N61" = -11627.479
  "H107" = 79200
  "J107" = 79200
  "L107" = 79200
  "N107" = -86880
  "H112" = 37137
  "J112" = 37137
  "L112" = 37137
  "N112" = -40091
  "H116" = 1480.6842
  "I116" = 1646
  "J116" = 860.75
  "K116" = 1646
  "L116" = 860.75
  "M116" = 648
  "N116" = -5448.75
  "H136" = 7582.74
  "I136" = 4498.407
  "J136" = 11203.479
  "K136" = 13495.221
  "L136" = 33610.437
  "M136" = -10945.221
  "N136" = -38710.437
  "H138" = 48047.11
  "J138" = 48047.11
  "L138" = 48047.11
  "N138" = -58327.11
}
foreach ($addr in $values.Keys) {
  $ws.Range($addr).Value = $values[$addr]
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$values = @{
  "H3" = 1480.6842
  "I3" = 1646
  "J3" = 860.75
  "K3" = 1646
  "L3" = 860.75
  "M3" = -1532
  "N3" = -1088.75
  "H107" = 2495.6956
  "I107" = 2033.3334
  "J107" = 3362.625
  "K107" = 2033.3334
  "L107" = 3362.625
  "M107" = -113.3334
  "N107" = -7202.625
  "H134" = 53492.25
  "I134" = 3957.3635
  "J134" = 114034.89
  "K134" = 11872.0905
  "L134" = 342104.67
  "M134" = -9337.0905
  "N134" = -347174.67
}
foreach ($addr in $values.Keys) {
  $ws.Range($addr).Value = $values[$addr]
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$values = @{
  "H16" = 2458.7144
  "I16" = 2142.2
  "K16" = 2142.2
  "M16" = -1855.2
  "H22" = 196.14285
  "I22" = 217.75
  "J22" = 167.33333
  "K22" = 217.75
  "L22" = 167.33333
  "M22" = 132.25
  "N22" = -867.3333299999999
  "H31" = 487244.3
  "I31" = 6199.6943
  "J31" = 942970.8
  "K31" = 6199.6943
  "L31" = 942970.8
  "M31" = -5904.6943
  "N31" = -943560.8
  "H34" = 487244.3
  "I34" = 6199.6943
  "J34" = 942970.8
  "K34" = 6199.6943
  "L34" = 942970.8
  "M34" = -5997.6943
  "N34" = -943374.8
  "H113" = 2458.7144
  "I113" = 2142.2
  "K113" = 2142.2
  "M113" = 27.80000000000018
  "H131" = 30326
  "J131" = 30326
  "L131" = 30326
  "N131" = -40406
  "H134" = 1800.9423
  "I134" = 1457.5
  "J134" = 4434
  "K134" = 4372.5
  "L134" = 13302
  "M134" = -1837.5
  "N134" = -18372
}
foreach ($addr in $values.Keys) {
  $ws.Range($addr).Value = $values[$addr]
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$values = @{
  "H68" = 1626.01
  "I68" = 1454.0339
  "J68" = 1873.4878
  "K68" = 4362.101699999999
  "L68" = 5620.463400000001
  "M68" = -3551.101699999999
  "N68" = -7242.463400000001
  "H71" = 1626.01
  "I71" = 1454.0339
  "J71" = 1873.4878
  "K71" = 13086.3051
  "L71" = 16861.3902
  "M71" = -9030.3051
  "N71" = -24973.3902
  "H120" = 11485.143
  "I120" = 8682.5
  "J120" = 15222
  "K120" = 26047.5
  "L120" = 45666
  "M120" = -21209.5
  "N120" = -55342
  "H131" = 1111.6111
  "J131" = 1439
  "L131" = 4317
  "N131" = -14397
}
foreach ($addr in $values.Keys) {
  $ws.Range($addr).Value = $values[$addr]
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$values = @{
  "H2" = 0
  "I2" = 0
  "K2" = 0
  "H18" = 30000004
  "I18" = 30000004
  "K18" = 30000004
  "M18" = -29999711
  "H107" = 473.875
  "I107" = 137.875
  "J107" = 809.875
  "K107" = 137.875
  "L107" = 809.875
  "M107" = 1782.125
  "N107" = -4649.875
  "H141" = 41000
  "J141" = 41000
  "L141" = 41000
  "N141" = -51360
}
foreach ($addr in $values.Keys) {
  $ws.Range($addr).Value = $values[$addr]
}
$ws.Range("M2").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$values = @{
  "H64" = 40150
  "J64" = 40150
  "L64" = 40150
  "N64" = -40600
  "H67" = 40150
  "J67" = 40150
  "L67" = 40150
  "N67" = -41710
}
foreach ($addr in $values.Keys) {
  $ws.Range($addr).Value = $values[$addr]
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$values = @{
  "H63" = 33149.4
  "I63" = 24500
  "J63" = 38915.668
  "K63" = 24500
  "L63" = 38915.668
  "M63" = -23876
  "N63" = -40163.668
  "H66" = 33149.4
  "I66" = 24500
  "J66" = 38915.668
  "K66" = 73500
  "L66" = 116747.004
  "M66" = -70380
  "N66" = -122987.004
  "H113" = 948
  "I113" = 609.25
  "J113" = 1173.8334
  "K113" = 1827.75
  "L113" = 3521.5002
  "M113" = 342.25
  "N113" = -7861.5002
  "H141" = 64949.5
  "J141" = 64949.5
  "L141" = 64949.5
  "N141" = -75309.5
}
foreach ($addr in $values.Keys) {
  $ws.Range($addr).Value = $values[$addr]
}
